# Update column G ("K" - strikeouts) values to regen save_data using K
# instead of Strike# (total strikes thrown). Values below are taken from
# the authoritative commit diff for xl/worksheets/sheet1.xml.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 3
    3  = 2
    4  = 0
    5  = 1
    6  = 0
    7  = 1
    8  = 1
    9  = 0
    10 = 0
    12 = 0
    13 = 0
    14 = 2
    16 = 0
    17 = 0
    18 = 2
    19 = 0
    20 = 1
    21 = 2
    22 = 1
    23 = 2
    24 = 1
    25 = 1
    26 = 0
    27 = 1
    28 = 2
    29 = 2
    30 = 1
    31 = 3
    32 = 0
    33 = 1
    34 = 0
    35 = 1
    36 = 0
    37 = 2
    38 = 0
    39 = 3
    40 = 0
    41 = 0
    42 = 0
    43 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
